# Avances Etiquetado Roboflow 6/4/2025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the data for the week of row 29 (D29 = 45753 -> 2/6/2025-4/6/2025)
$ws.Range("E29").Value = 127
$ws.Range("F29").Value = 234
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 1012
$ws.Range("J29").Value = "N/A"

# Update the view: scroll/zoom and selection as left by the author when saving
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I32").Select()
